$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the style of the last existing header cell (E1) onto the two new
# header cells so they pick up the same cellXfs entry (center/center,
# no new style record) instead of synthesizing a brand-new style.
$ws.Range("E1").Copy($ws.Range("F1:G1")) | Out-Null

# New header labels
$ws.Range("F1").Value = "BL"
$ws.Range("G1").Value = "Operantig Freq"

# New data row values
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 0

# Match the selection shown in the edited file
$ws.Range("G2").Select()
